$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.022877812385559
$ws.Range("B1").Value = 1.539020776748657
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 0.3124167025089264
